# Adds a new reservation record (row 21) to the export sheet, mirroring the
# row 19/20 layout, data types and formatting already used in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21 : new reservation ------------------------------------------------
# General-purpose ("text"/"number", 10pt Calibri, no border) cells first, so
# the engine allocates the shared cellXf (General / font 5) once and reuses it
# for every plain cell below.

$ws.Range("A21").Value = 11178448
$ws.Range("A21").Font.Size = 10

$ws.Range("B21").Value = "PH11178448"
$ws.Range("B21").Font.Size = 10

$ws.Range("C21").Value = "ACSI"
$ws.Range("C21").Font.Size = 10

$ws.Range("D21").Value = "Kamp Koren"
$ws.Range("D21").Font.Size = 10

$ws.Range("E21").Value = "74958700e8681600f4599f30aadcc810"
$ws.Range("E21").Font.Size = 10

# Arrival / Departure - dd/mm/yy date cells (distinct number format). Values
# are supplied as the raw 1900-date-system serials (17/20 Sep 2021) so the
# engine doesn't pre-assign its own default date numFmt before we set ours.
$ws.Range("F21").Value = 44456
$ws.Range("F21").Font.Size = 10
$ws.Range("F21").NumberFormat = "dd/mm/yy"

$ws.Range("G21").Value = 44459
$ws.Range("G21").Font.Size = 10
$ws.Range("G21").NumberFormat = "dd/mm/yy"

$ws.Range("H21").Value = "Camper / Caravan site"
$ws.Range("H21").Font.Size = 10

$ws.Range("I21").Value = "Duin Jurjen, mr"
$ws.Range("I21").Font.Size = 10

$ws.Range("J21").Value = 101
$ws.Range("J21").Font.Size = 10

$ws.Range("K21").Value = 12
$ws.Range("K21").Font.Size = 10

$ws.Range("L21").Value = "EUR"
$ws.Range("L21").Font.Size = 10

$ws.Range("M21").Value = "cc"
$ws.Range("M21").Font.Size = 10

$ws.Range("N21").Value = "OK"
$ws.Range("N21").Font.Size = 10

# Date of creation - dd/mm/yy date cell.
$ws.Range("O21").Value = 44454
$ws.Range("O21").Font.Size = 10
$ws.Range("O21").NumberFormat = "dd/mm/yy"

$ws.Range("P21").Value = "marijkeduin@kpnmail.nl"
$ws.Range("P21").Font.Size = 10

$ws.Range("Q21").Value = 641122652
$ws.Range("Q21").Font.Size = 10

$ws.Range("R21").Value = "Netherlands"
$ws.Range("R21").Font.Size = 10

$ws.Range("S21").Value = "NL"
$ws.Range("S21").Font.Size = 10

$ws.Range("T21").Value = "Vaart noordzijde, 33"
$ws.Range("T21").Font.Size = 10

$ws.Range("U21").Value = "8426AR"
$ws.Range("U21").Font.Size = 10

$ws.Range("V21").Value = "Appelscha"
$ws.Range("V21").Font.Size = 10

$ws.Range("W21").Value = "#ACSI"
$ws.Range("W21").Font.Size = 10

# Column X (B2B discount applied) is left blank for this row, same as the
# source data.

$ws.Range("Y21").Value = "Jurjen Duin (1954-06-26), Marijke Duin (1968-07-20)"
$ws.Range("Y21").Font.Size = 10

$ws.Range("Z21").Value = 0
$ws.Range("Z21").Font.Size = 10

$ws.Range("AA21").Value = "Best available rate"
$ws.Range("AA21").Font.Size = 10

# Match the slightly reduced row height used for this data row.
$ws.Rows.Item(21).RowHeight = 12.8

# Leave the view scrolled/selected the way it was after entering this row
# (top-left visible cell around row 15, active cell D16).
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1
$ws.Range("D16").Select()
